# Auto_Fernando_Fuentes_11B.docx edits:
#   1. Fernando José Fuentes Castillo's autoevaluación score: 10 -> 9.5
#   2. Mandhy Guadalupe Masin Rodríguez's coevaluación score: 10 -> 9
#   3. Avril Fernanda Paz Pinto's coevaluación score: 9.5 -> 9
#      (the ".5" lived in its own run; that run is removed entirely)
#   4. Avril's justification paragraph gets a new trailing sentence,
#      added as its own new run.
#
# Note: this runtime's Range.Find.Execute(..., Replace:=wdReplaceAll/One)
# operates over the *whole* document content regardless of which Range
# object it is invoked on, so it is avoided here in favour of directly
# editing Range.Text / Range.Delete() / Range.InsertAfter() on precisely
# bounded Range objects (each Table is re-fetched fresh after a mutation
# since prior Table/Cell handles can go stale once the content shifts).

$d = $word.ActiveDocument

# --- 1. Fernando's score cell: "10" -> "9.5" ------------------------------
$t = $d.Tables(1)
$cell = $t.Cell(4, 4)
$r = $cell.Range
$r.End = $r.End - 1          # drop trailing cell-end mark
$r.Text = "9.5"

# --- 2. Mandhy's score cell: "10" -> "9" ----------------------------------
$t = $word.ActiveDocument.Tables(1)
$cell = $t.Cell(9, 4)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "9"

# --- 3. Avril's score cell: drop the ".5" run, leaving just "9" ----------
$t = $word.ActiveDocument.Tables(1)
$cell = $t.Cell(11, 4)
$r = $cell.Range
$tail = $d.Range($r.Start + 1, $r.End - 1)   # the ".5" portion only
$tail.Delete()

# --- 4. Append a new sentence/run to Avril's justification paragraph -----
$t = $word.ActiveDocument.Tables(1)
$cell = $t.Cell(12, 1)
$para = $cell.Range.Paragraphs(1)
$r = $para.Range
$r.End = $r.End - 1          # drop the paragraph mark
$r.Collapse(0)
$r.InsertAfter(" Mostró interés en todo momento por el trabajo.")

"done"
